$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Stash the current cell formats (style indices) of A1 (SampleID style)
#        and B1 (#Chr/Start/Stop/Ref/Call style) into scratch cells far to the
#        right, so we can re-apply them after the columns are rebuilt.
$ws.Range("A1").Copy()
$ws.Range("AA1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B1").Copy()
$ws.Range("AB1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# --- 2. Drop columns A:F completely (this removes their contents AND the
#        <cols> hidden/width overrides), then insert 6 fresh, fully default
#        columns back in their place so no hidden/width formatting remains.
$ws.Range("A1:F1").EntireColumn.Delete()
$ws.Range("A1:F1").EntireColumn.Insert()

# --- 3. Write the new header order: #Chr, Start, Stop, Ref, Call, SampleID
$ws.Range("A1").Value2 = "#Chr"
$ws.Range("B1").Value2 = "Start"
$ws.Range("C1").Value2 = "Stop"
$ws.Range("D1").Value2 = "Ref"
$ws.Range("E1").Value2 = "Call"
$ws.Range("F1").Value2 = "SampleID"

# --- 4. Re-apply the stashed formats: A1:E1 get the old B1 style, F1 gets
#        the old A1 style.
$ws.Range("AB1").Copy()
$ws.Range("A1:E1").PasteSpecial(-4122) | Out-Null
$ws.Range("AA1").Copy()
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# --- 5. Clean up the scratch cells used to stash formats.
$ws.Range("AA1:AB1").Clear() | Out-Null

# --- 6. Rebuild the AutoFilter so it only covers F1:Y1.
$ws.AutoFilterMode = $false
$ws.Range("F1:Y1").AutoFilter() | Out-Null

# --- 7. Update the worksheet-scoped _FilterDatabase defined name to match.
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -like "*_FilterDatabase") {
        $n.RefersTo = "=Tier3!`$F`$1:`$Y`$1"
    }
}

# --- 8. Update the visible selection to a single cell, C1.
$ws.Range("C1").Select() | Out-Null

Write-Host "done"
